$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 updates
$ws.Range("M4").Value = 1.1
$ws.Range("N4").Value = 7

# Row 5 updates
$ws.Range("G5").Value = 2.3
$ws.Range("J5").Value = 3.1
$ws.Range("M5").Value = 1.06
$ws.Range("N5").Value = 9.5
$ws.Range("O5").Value = 1.33
$ws.Range("P5").Value = 3.25
$ws.Range("Q5").Value = 2.08
$ws.Range("R5").Value = 1.73
$ws.Range("Y5").Value = 9.5
$ws.Range("AB5").Value = 29
$ws.Range("AC5").Value = 9.5
$ws.Range("AH5").Value = 8.5
$ws.Range("AN5").Value = 4.33
$ws.Range("AO5").Value = 13
$ws.Range("AP5").Value = 23
$ws.Range("AQ5").Value = 41
